# Add a "validate" step (pages.LoginPage) between "login" and "logout" on the
# LoginPage sheet, and clean out the now-empty placeholder cells on the
# AddEmployeePage sheet (still "in progress" - only openBrowser/launchUrl/login
# are filled in for rows 3-8). Also move the active/selected tab from
# AddEmployeePage back to LoginPage.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # LoginPage
$ws2 = $wb.Worksheets.Item(2)   # AddEmployeePage

# --- LoginPage: insert a new column E holding the "validate" step ---------
# This shifts the old E (logout) -> F and old F (closeBrowser) -> G.
$ws1.Columns("E:E").Insert()

# Give the newly inserted column a sensible width (matches column D's
# "pages.X" style width).
$ws1.Columns("E:E").ColumnWidth = 14.75

# Header row
$ws1.Range("E1").Value = "pages.LoginPage"

# Row 2 keeps the full login/validate/logout/closeBrowser sequence.
$ws1.Range("E2").Value = "validate"

# Rows 3-8 only use validate + closeBrowser (no explicit logout step, same
# pattern the sheet already used for the closeBrowser-only rows).
foreach ($r in 3..8) {
    $ws1.Range("E$r").Value = "validate"
    $ws1.Range("F$r").ClearContents()
}

# --- AddEmployeePage: drop the leftover empty E column for rows 3-8 -------
$ws2.Range("E3:E8").ClearContents()

# --- Tab selection: move the active sheet back to LoginPage ---------------
$null = $ws1.Select()
$null = $ws1.Range("F3").Select()
